# Weekly fruit/vegetable price log: insert a new weekly record at row 195
# ("Hortaliza, Femacal de La Calera - Espinaca"), which pushes all the
# existing records from row 195 down to row 196 (through row 287).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting the entire row shifts rows 195:286 down to 196:287 and grows
# the used range from A1:R286 to A1:R287 automatically, carrying each
# cell's formatting (e.g. the date style on column D) along with it.
$ws.Rows.Item(195).Insert()

# Populate the newly-opened row 195 with this week's record.
$ws.Range("A195").Value = 3
$ws.Range("B195").Value = "Femacal de La Calera"
$ws.Range("C195").Value = "Coquimbo"
$ws.Range("D195").Value = 44609
$ws.Range("E195").Value = 5
$ws.Range("F195").Value = 100112012
$ws.Range("G195").Value = "Espinaca"
$ws.Range("H195").Value = "Sin especificar"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 125
$ws.Range("K195").Value = 4500
$ws.Range("L195").Value = 5000
$ws.Range("M195").Value = 4740
$ws.Range("N195").Value = "$/docena de atados (3 kilos)"
$ws.Range("O195").Value = "Provincia de Quillota"
$ws.Range("P195").Value = 1580
$ws.Range("Q195").Value = 3
$ws.Range("R195").Value = "Hortaliza"
